$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update B2 (CasesTab query) and C2 (StatQuery) to their new text ---
$ws.Range("B2").Value = 'MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
  WHERE c.ethnicity = "HISPANIC_OR_LATINO"
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity'
$ws.Range("C2").Value = 'MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
       WHERE c.ethnicity = "HISPANIC_OR_LATINO"
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files'

# --- Row 3: new FilesTab row ---
$ws.Range("A3").Value = 'FilesTab'
$ws.Range("B3").Value = 'MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
 WHERE c.ethnicity = "HISPANIC_OR_LATINO"
WITH
    f, parent, c, a, ct,
    [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`'
$ws.Range("C3").Value = 'MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
        WHERE c.ethnicity = "HISPANIC_OR_LATINO"
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files'
$ws.Range("D3").Value = $ws.Range("D2").Value()
$ws.Range("E3").Value = $ws.Range("E2").Value()

# --- Apply wrap-text style to B3/C3, matching B2/C2 ---
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 188.5
$ws.Rows.Item(3).RowHeight = 409.5

# --- View: selection + zoom ---
$ws.Range("C3").Select()
$excel.ActiveWindow.Zoom = 70
